$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format first so dotted numeric-looking
# strings (e.g. "1.003", "24.944.83") are stored verbatim as text, matching
# the original inline-string cells, instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.944.83"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "1.708.27"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "313.13"
$ws.Range("E5").Value = "  +2.47%  "
$ws.Range("D6").Value = "0.9980"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  +1.55%  "
$ws.Range("D8").Value = "49.48"
$ws.Range("E8").Value = "  +3.66%  "
$ws.Range("D9").Value = "0.3442"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").Value = "1.226"
$ws.Range("E10").Value = "  +6.09%  "
$ws.Range("D11").Value = "0.07539"
$ws.Range("E11").Value = "  +4.61%  "
$ws.Range("D12").Value = "0.9990"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("E13").Value = "  +6.14%  "
$ws.Range("D14").Value = "6.326"
$ws.Range("E14").Value = "  +4.16%  "
$ws.Range("D15").Value = "7.078"
$ws.Range("E15").Value = "  +5.69%  "
$ws.Range("D16").Value = "1.705.69"
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("E17").Value = "  +3.07%  "
$ws.Range("D18").Value = "0.06729"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").Value = "0.9970"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "84.16"
$ws.Range("E20").Value = "  +4.41%  "
$ws.Range("D21").Value = "17.33"
$ws.Range("E21").Value = "  +6.19%  "
$ws.Range("D22").Value = "6.395"
$ws.Range("E22").Value = "  +5.52%  "
$ws.Range("E23").Value = "  +8.67%  "
$ws.Range("D24").Value = "24.935.21"
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("D25").Value = "2.441"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  +6.51%  "
$ws.Range("D27").Value = "20.42"
$ws.Range("E27").Value = "  +5.84%  "
$ws.Range("D28").Value = "149.86"
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("D29").Value = "132.90"
$ws.Range("E29").Value = "  +4.78%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "1.253"
$ws.Range("E30").Value = "  +30.39%  "
$ws.Range("B31").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C31").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D31").Value = "1.895.00"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("D32").Value = "6.835"
$ws.Range("E32").Value = "  +9.68%  "
$ws.Range("D33").Value = "4.224"
$ws.Range("E33").Value = "  +4.03%  "
$ws.Range("D34").Value = "13.91"
$ws.Range("E34").Value = "  +14.07%  "
$ws.Range("D35").Value = "0.08787"
$ws.Range("E35").Value = "  +4.60%  "
$ws.Range("D36").Value = "1.773"
$ws.Range("E36").Value = "  +4.54%  "
$ws.Range("D37").Value = "5.638"
$ws.Range("E37").Value = "  +6.77%  "
$ws.Range("D38").Value = "0.06671"
$ws.Range("E38").Value = "  +4.06%  "
$ws.Range("D39").Value = "9.187"
$ws.Range("E39").Value = "  +4.71%  "
$ws.Range("D40").Value = "0.02417"
$ws.Range("E40").Value = "  +5.43%  "
$ws.Range("D41").Value = "0.2239"
$ws.Range("E41").Value = "  +7.74%  "
$ws.Range("D42").Value = "1.272"
$ws.Range("E42").Value = "  +3.34%  "
$ws.Range("D43").Value = "0.6480"
$ws.Range("E43").Value = "  +6.73%  "
$ws.Range("D44").Value = "0.9974"
$ws.Range("D45").Value = "13.89"
$ws.Range("E45").Value = "  +7.50%  "
$ws.Range("D46").Value = "0.6170"
$ws.Range("E46").Value = "  +5.32%  "
$ws.Range("D47").Value = "3.843"
$ws.Range("E47").Value = "  +2.40%  "
$ws.Range("D48").Value = "2.126"
$ws.Range("E48").Value = "  +6.23%  "
$ws.Range("D49").Value = "129.47"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("D50").Value = "0.07333"
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("D51").Value = "80.40"
$ws.Range("E51").Value = "  +6.89%  "

# Reset column D back to the default (General/Normal) style so no residual
# explicit cell style reference is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"

